# Insert a new weekly price record at row 83, pushing the existing
# rows 83-148 down to 84-149 (matches the "Fruta / hortaliza, semanal"
# weekly-refresh commit: dimension grows from A1:R148 to A1:R149).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 83..148 down to 84..149, leaving a blank row 83.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new observation.
$ws.Range("A83").Value2 = 9
$ws.Range("B83").Value2 = 'Vega Central Mapocho de Santiago'
$ws.Range("C83").Value2 = 'Metropolitana'
$ws.Range("D83").Value2 = 45236
$ws.Range("E83").Value2 = 13
$ws.Range("F83").Value2 = 100114007
$ws.Range("G83").Value2 = 'Jengibre'
$ws.Range("H83").Value2 = 'Sin especificar'
$ws.Range("I83").Value2 = 'Primera'
$ws.Range("J83").Value2 = 520
$ws.Range("K83").Value2 = 23000
$ws.Range("L83").Value2 = 25000
$ws.Range("M83").Value2 = 24000
$ws.Range("N83").Value2 = '$/caja 13 kilos'
$ws.Range("O83").Value2 = 'Perú'
$ws.Range("P83").Value2 = 1846
$ws.Range("Q83").Value2 = 13
$ws.Range("R83").Value2 = 'Hortaliza'
